# Apply "Edited proto, added authentication server responses to chat server"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11 header: D/E relabelled, F becomes the new "Server " response column ---
$ws.Range("D11").Value = "Server"
$ws.Range("E11").Value = "Authentication"

# --- Rows 12-15: new "Server " response column (F) ---
$ws.Range("F12").Value = "Register Success"
$ws.Range("F14").Value = "Login Success"
$ws.Range("F13").Value = "Register Failure"
$ws.Range("F15").Value = "Login Failure"

# --- Row 10: add RequestId note ---
$ws.Range("D10").Value = "RequestId=client index"

$ws.Range("F11").Value = "Server "

# --- New rows 17-20 ---
$ws.Range("A17").Value = 5
$ws.Range("D17").Value = "Register/Login Success"

$ws.Range("A18").Value = 6
$ws.Range("D18").Value = "Register/Login Failure"

$ws.Range("A19").Value = 7

$ws.Range("A20").Value = 8

# --- View state: scroll + selection to match author's final position ---
$win = $excel.ActiveWindow
$win.ScrollRow = 5
[void]$ws.Range("C18").Select()

# --- Column widths E:F recompute (as Excel's bestFit would after the edits) ---
$ws.Columns.Item(5).ColumnWidth = 12.25
$ws.Columns.Item(6).ColumnWidth = 13.5
